$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.531238913048924
$ws.Range("C2").Value = 0.2923111209445892
$ws.Range("D2").Value = 0.05423374491933686
$ws.Range("F2").Value = 4.845508045540868
$ws.Range("G2").Value = 0.002564576935133909
$ws.Range("J2").Value = 0.2314661405696157
$ws.Range("M2").Value = 0.5003697753631613
$ws.Range("B3").Value = 1.460752453487544
$ws.Range("C3").Value = 0.2756330872252875
$ws.Range("D3").Value = 0.05186890427511059
$ws.Range("F3").Value = 4.659147865040865
$ws.Range("G3").Value = 0.002571257130335096
$ws.Range("J3").Value = 0.2290218502733552
$ws.Range("M3").Value = 0.4894542456953488
$ws.Range("B4").Value = 1.418927386970267
$ws.Range("C4").Value = 0.2657263352744792
$ws.Range("D4").Value = 0.05043720563763543
$ws.Range("F4").Value = 4.546096390128184
$ws.Range("G4").Value = 0.002575568007025162
$ws.Range("J4").Value = 0.227594385278465
$ws.Range("M4").Value = 0.4832266224057946
$ws.Range("B5").Value = 1.402246489376921
$ws.Range("C5").Value = 0.2617723774270644
$ws.Range("D5").Value = 0.04985870104599144
$ws.Range("F5").Value = 4.500364531567527
$ws.Range("G5").Value = 0.002577377537643575
$ws.Range("J5").Value = 0.2270310043539823
$ws.Range("M5").Value = 0.4808077124416741
$ws.Range("B6").Value = 1.399498502753232
$ws.Range("C6").Value = 0.261120824081388
$ws.Range("D6").Value = 0.04976293278848942
$ws.Range("F6").Value = 4.49279094311882
$ws.Range("G6").Value = 0.00257768120465452
$ws.Range("J6").Value = 0.2269385585848909
$ws.Range("M6").Value = 0.4804132230009444
$ws.Range("B7").Value = 1.418700955312659
$ws.Range("C7").Value = 0.2656726752521763
$ws.Range("D7").Value = 0.05042938402936414
$ws.Range("F7").Value = 4.545478278234725
$ws.Range("G7").Value = 0.002575592196911746
$ws.Range("J7").Value = 0.2275867132940377
$ws.Range("M7").Value = 0.4831935192887471
$ws.Range("B8").Value = 1.506631909905877
$ws.Range("C8").Value = 0.2864908157531261
$ws.Range("D8").Value = 0.05341401680580304
$ws.Range("F8").Value = 4.780960851947469
$ws.Range("G8").Value = 0.002566836969119525
$ws.Range("J8").Value = 0.230608054482957
$ws.Range("M8").Value = 0.4965073353848197
$ws.Range("B9").Value = 1.690724654860333
$ws.Range("C9").Value = 0.3300004243485546
$ws.Range("D9").Value = 0.05943681045265237
$ws.Range("F9").Value = 5.254037188587375
$ws.Range("G9").Value = 0.002551318568239354
$ws.Range("J9").Value = 0.2371206965491623
$ws.Range("M9").Value = 0.5264037582998
$ws.Range("B10").Value = 1.833279031059362
$ws.Range("C10").Value = 0.3636616835578366
$ws.Range("D10").Value = 0.06397779582782448
$ws.Range("F10").Value = 5.609070592074772
$ws.Range("G10").Value = 0.00254091020905592
$ws.Range("J10").Value = 0.2422727801936304
$ws.Range("M10").Value = 0.5507139058930619
$ws.Range("B11").Value = 1.899760247329311
$ws.Range("C11").Value = 0.3793561631404145
$ws.Range("D11").Value = 0.06607153650551822
$ws.Range("F11").Value = 5.77233383201343
$ws.Range("G11").Value = 0.002536387960333281
$ws.Range("J11").Value = 0.2446983468357686
$ws.Range("M11").Value = 0.5622908031917859
$ws.Range("B12").Value = 1.92517306475861
$ws.Range("C12").Value = 0.3853551785644811
$ws.Range("D12").Value = 0.06686862834422413
$ws.Range("F12").Value = 5.834419894059295
$ws.Range("G12").Value = 0.002534705849767961
$ws.Range("J12").Value = 0.2456287710267375
$ws.Range("M12").Value = 0.5667498107430617
$ws.Range("B13").Value = 1.919689332688108
$ws.Range("C13").Value = 0.3840606835469771
$ws.Range("D13").Value = 0.06669676854590989
$ws.Range("F13").Value = 5.821036757350953
$ws.Range("G13").Value = 0.002535066774849561
$ws.Range("J13").Value = 0.2454278552978408
$ws.Range("M13").Value = 0.5657861345498674
$ws.Range("B14").Value = 1.901846190861818
$ws.Range("C14").Value = 0.3798485807418785
$ws.Range("D14").Value = 0.06613702760093076
$ws.Range("F14").Value = 5.777436389309742
$ws.Range("G14").Value = 0.002536248964610451
$ws.Range("J14").Value = 0.2447746537701789
$ws.Range("M14").Value = 0.5626561401144272
$ws.Range("B15").Value = 1.890947815624656
$ws.Range("C15").Value = 0.3772758490956676
$ws.Range("D15").Value = 0.06579472819669263
$ws.Range("F15").Value = 5.75076427663322
$ws.Range("G15").Value = 0.002536977038639102
$ws.Range("J15").Value = 0.2443761049442514
$ws.Range("M15").Value = 0.5607487245929548
$ws.Range("B16").Value = 1.828967437949018
$ws.Range("C16").Value = 0.3626437761931527
$ws.Range("D16").Value = 0.06384154657409624
$ws.Range("F16").Value = 5.598437088680043
$ws.Range("G16").Value = 0.002541210009088209
$ws.Range("J16").Value = 0.2421159231099566
$ws.Range("M16").Value = 0.5499678034071707
$ws.Range("B17").Value = 1.791364765589947
$ws.Range("C17").Value = 0.3537659326930793
$ws.Range("D17").Value = 0.06265065759259869
$ws.Range("F17").Value = 5.505445738596109
$ws.Range("G17").Value = 0.002543861100233634
$ws.Range("J17").Value = 0.2407504417807189
$ws.Range("M17").Value = 0.5434871634688534
$ws.Range("B18").Value = 1.769890173959595
$ws.Range("C18").Value = 0.3486955104662854
$ws.Range("D18").Value = 0.06196831488536247
$ws.Range("F18").Value = 5.45212436975342
$ws.Range("G18").Value = 0.002545405957478441
$ws.Range("J18").Value = 0.2399727491121126
$ws.Range("M18").Value = 0.5398083993611209
$ws.Range("B19").Value = 1.762645512940253
$ws.Range("C19").Value = 0.34698488617164
$ws.Range("D19").Value = 0.0617377299458326
$ws.Range("F19").Value = 5.434098696556475
$ws.Range("G19").Value = 0.002545932464063515
$ws.Range("J19").Value = 0.2397107532202938
$ws.Range("M19").Value = 0.5385711846811603
$ws.Range("B20").Value = 1.795351728650928
$ws.Range("C20").Value = 0.3547072741728812
$ws.Range("D20").Value = 0.06277715656759142
$ws.Range("F20").Value = 5.515327704641777
$ws.Range("G20").Value = 0.002543576816664795
$ws.Range("J20").Value = 0.2408950020539322
$ws.Range("M20").Value = 0.5441719912477438
$ws.Range("B21").Value = 1.90708067474975
$ws.Range("C21").Value = 0.3810842537725136
$ws.Range("D21").Value = 0.06630132040594106
$ws.Range("F21").Value = 5.790235694641751
$ws.Range("G21").Value = 0.002535900904272725
$ws.Range("J21").Value = 0.2449661903696025
$ws.Range("M21").Value = 0.5635734527024994
$ws.Range("B22").Value = 1.981489708886841
$ws.Range("C22").Value = 0.3986491717324725
$ws.Range("D22").Value = 0.0686293639955835
$ws.Range("F22").Value = 5.971434399163229
$ws.Range("G22").Value = 0.002531061161089522
$ws.Range("J22").Value = 0.2476964645015727
$ws.Range("M22").Value = 0.5766913866328593
$ws.Range("B23").Value = 1.941648193109984
$ws.Range("C23").Value = 0.3892442899132789
$ws.Range("D23").Value = 0.06738450486254521
$ws.Range("F23").Value = 5.87458212743519
$ws.Range("G23").Value = 0.002533628101875833
$ws.Range("J23").Value = 0.2462328566806562
$ws.Range("M23").Value = 0.5696498255602762
$ws.Range("B24").Value = 1.793548775332965
$ws.Range("C24").Value = 0.3542815892044757
$ws.Range("D24").Value = 0.06271995918521611
$ws.Range("F24").Value = 5.510859630557491
$ws.Range("G24").Value = 0.002543705276749895
$ws.Range("J24").Value = 0.2408296235004741
$ws.Range("M24").Value = 0.5438622341272321
$ws.Range("B25").Value = 1.639654726583274
$ws.Range("C25").Value = 0.3179365654974617
$ws.Range("D25").Value = 0.05778801400525424
$ws.Range("F25").Value = 5.124788945859109
$ws.Range("G25").Value = 0.002555341369296549
$ws.Range("J25").Value = 0.2352950158576377
$ws.Range("M25").Value = 0.5179068035163965
